$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> column -> new value, derived from the authoritative OOXML diff.
$rowData = [ordered]@{
    2 = @{ "D" = "63.457.65"; "E" = "  +0.01%  " }
    3 = @{ "D" = "3.084.02"; "E" = "  -0.23%  " }
    4 = @{ "E" = "  -0.01%  " }
    5 = @{ "D" = "547.37"; "E" = "  -0.23%  " }
    6 = @{ "D" = "139.39"; "E" = "  +1.92%  " }
    7 = @{ "D" = "1.00"; "E" = "  +0.07%  " }
    8 = @{ "D" = "3.078.63"; "E" = "  -0.53%  " }
    9 = @{ "E" = "  +0.43%  " }
    10 = @{ "E" = "  +1.19%  " }
    11 = @{ "D" = "6.39"; "E" = "  +0.82%  " }
    12 = @{ "E" = "  -2.43%  " }
    13 = @{ "E" = "  +3.81%  " }
    14 = @{ "D" = "34.98"; "E" = "  -1.17%  " }
    15 = @{ "D" = "3.586.23"; "E" = "  -0.13%  " }
    16 = @{ "D" = "63.489.62"; "E" = "  +0.06%  " }
    17 = @{ "E" = "  +1.09%  " }
    18 = @{ "D" = "3.083.08"; "E" = "  -0.33%  " }
    19 = @{ "E" = "  -1.16%  " }
    20 = @{ "D" = "475.50"; "E" = "  -2.58%  " }
    21 = @{ "E" = "  -0.37%  " }
    22 = @{ "D" = "0.703"; "E" = "  -1.97%  " }
    23 = @{ "D" = "7.09"; "E" = "  -2.51%  " }
    24 = @{ "D" = "78.71"; "E" = "  -0.42%  " }
    25 = @{ "D" = "12.24"; "E" = "  -0.89%  " }
    26 = @{ "D" = "0.999"; "E" = "  -0.06%  " }
    27 = @{ "E" = "  -1.05%  " }
    28 = @{ "D" = "7.95"; "E" = "  -6.47%  " }
    30 = @{ "D" = "26.28"; "E" = "  -1.20%  " }
    31 = @{ "B" = "Mantle"; "C" = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; "D" = "1.17"; "E" = "  +3.83%  " }
    32 = @{ "B" = "ImmutableX"; "C" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; "D" = "1.90"; "E" = "  -3.48%  " }
    33 = @{ "D" = "58.73"; "E" = "  -0.18%  " }
    34 = @{ "E" = "  -7.31%  " }
    35 = @{ "D" = "5.53"; "E" = "  +8.57%  " }
    36 = @{ "D" = "6.03"; "E" = "  -0.41%  " }
    37 = @{ "D" = "489.17"; "E" = "  -3.40%  " }
    38 = @{ "D" = "3.259.41"; "E" = "  +3.55%  " }
    39 = @{ "D" = "0.0403"; "E" = "  +1.00%  " }
    40 = @{ "D" = "0.0798"; "E" = "  -0.55%  " }
    41 = @{ "E" = "  -0.91%  " }
    42 = @{ "D" = "8.16"; "E" = "  +0.05%  " }
    43 = @{ "D" = "2.60"; "E" = "  -1.10%  " }
    44 = @{ "D" = "0.253"; "E" = "  -0.89%  " }
    45 = @{ "E" = "  +0.07%  " }
    46 = @{ "B" = "Monero"; "C" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; "D" = "124.71"; "E" = "  +3.75%  " }
    47 = @{ "B" = "InjectiveProtocol"; "C" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; "D" = "25.46"; "E" = "  +0.08%  " }
    48 = @{ "E" = "  -1.82%  " }
    49 = @{ "D" = "0.0₃0531"; "E" = "  +5.02%  " }
    50 = @{ "E" = "  +0.44%  " }
    51 = @{ "E" = "  -0.35%  " }
}

foreach ($rowKey in $rowData.Keys) {
    $row = [int]$rowKey
    $cols = $rowData[$rowKey]
    foreach ($col in $cols.Keys) {
        $newValue = $cols[$col]
        $cell = $ws.Range("$col$row")
        # Force text storage so numeric-looking strings (e.g. "1.00") are not
        # reinterpreted as numbers by the COM Value setter, then restore the
        # default (unstyled) cell style so no stray formatting is introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    }
}

Write-Output "Updated $($rowData.Count) rows"
